# Apply the edit described by the diff:
#  1. Move the "active/selected" tab from PP (4th sheet) to DMD (5th sheet),
#     and change the sheetView of DMD to be the selected one, with a new
#     zoom level and selection.
#  2. On the DMD sheet, delete the range F5:F8 and shift the remaining
#     cells (G:K) left by one column.

$wb = $excel.ActiveWorkbook

$wsPP  = $wb.Worksheets.Item("PP")
$wsDMD = $wb.Worksheets.Item("DMD")

# --- Sheet DMD: delete the F5:K8 block, shifting cells left (within rows 5-8 only) ---
$wsDMD.Range("F5:K8").Delete(-4159)  # xlShiftToLeft = -4159

# --- Update the view/selection state ---
# Previously PP's sheetView had zoomScale/zoomScaleNormal = 130 and was the
# selected ("tabSelected") sheet with selection E9. Now DMD becomes the
# selected sheet, zoom 120, with selection K13.
$wsDMD.Select()
$excel.ActiveWindow.Zoom = 120
$wsDMD.Range("K13").Select()

$wb.Save()
